$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3418.8484
$ws.Range("J112").Value = 3458.9688
$ws.Range("L112").Value = 10376.9064
$ws.Range("N112").Value = -12592.9064

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 13898986
$ws.Range("I113").Value = 35725000
$ws.Range("K113").Value = 35725000
$ws.Range("M113").Value = -35721746

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4643.337
$ws.Range("I138").Value = 1333.409
$ws.Range("J138").Value = 5601.4736
$ws.Range("K138").Value = 4000.227
$ws.Range("L138").Value = 16804.4208
$ws.Range("M138").Value = 1139.773
$ws.Range("N138").Value = -27084.4208

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4940.421
$ws.Range("I141").Value = 5725.1665
$ws.Range("K141").Value = 17175.4995
$ws.Range("M141").Value = -11995.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2594.68
$ws.Range("I32").Value = 2060.0857
$ws.Range("J32").Value = 10079
$ws.Range("K32").Value = 2060.0857
$ws.Range("L32").Value = 10079
$ws.Range("M32").Value = -1773.0857
$ws.Range("N32").Value = -10653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1456.8235
$ws.Range("I45").Value = 1429.1818
$ws.Range("K45").Value = 1429.1818
$ws.Range("M45").Value = -1052.1818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2839.6875
$ws.Range("I61").Value = 2719.6897
$ws.Range("J61").Value = 3999.6667
$ws.Range("K61").Value = 2719.6897
$ws.Range("L61").Value = 3999.6667
$ws.Range("M61").Value = -2507.6897
$ws.Range("N61").Value = -4423.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2131.257
$ws.Range("I74").Value = 1730.8518
$ws.Range("K74").Value = 1730.8518
$ws.Range("M74").Value = -856.8517999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2131.257
$ws.Range("I77").Value = 1730.8518
$ws.Range("K77").Value = 8654.259
$ws.Range("M77").Value = -4286.259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 627.5
$ws.Range("I97").Value = 666.3333
$ws.Range("J97").Value = 511
$ws.Range("K97").Value = 666.3333
$ws.Range("L97").Value = 511
$ws.Range("M97").Value = -170.3333
$ws.Range("N97").Value = -1503

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 10527636
$ws.Range("I110").Value = 13334685
$ws.Range("J110").Value = 1203.5
$ws.Range("K110").Value = 13334685
$ws.Range("L110").Value = 1203.5
$ws.Range("M110").Value = -13332640
$ws.Range("N110").Value = -5293.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 12505.059
$ws.Range("I122").Value = 4762.4546
$ws.Range("J122").Value = 26699.834
$ws.Range("K122").Value = 14287.3638
$ws.Range("L122").Value = 80099.50199999999
$ws.Range("M122").Value = -11837.3638
$ws.Range("N122").Value = -84999.50199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 75594.60000000001
$ws.Range("J123").Value = 75594.60000000001
$ws.Range("L123").Value = 75594.60000000001
$ws.Range("N123").Value = -85394.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2953.23
$ws.Range("I132").Value = 2939.2625
$ws.Range("J132").Value = 3009.1
$ws.Range("K132").Value = 8817.787499999999
$ws.Range("L132").Value = 9027.299999999999
$ws.Range("M132").Value = -6287.787499999999
$ws.Range("N132").Value = -14087.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2839.6875
$ws.Range("I136").Value = 2719.6897
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 8159.0691
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = -5609.0691
$ws.Range("N136").Value = -17099.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2041.1364
$ws.Range("I20").Value = 1850.0667
$ws.Range("J20").Value = 2450.5715
$ws.Range("K20").Value = 1850.0667
$ws.Range("L20").Value = 2450.5715
$ws.Range("M20").Value = -1603.0667
$ws.Range("N20").Value = -2944.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13842
$ws.Range("I82").Value = 13842
$ws.Range("K82").Value = 13842
$ws.Range("M82").Value = -13459

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 13842
$ws.Range("I85").Value = 13842
$ws.Range("K85").Value = 13842
$ws.Range("M85").Value = -12516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2378.1794
$ws.Range("I86").Value = 2280.4546
$ws.Range("J86").Value = 2915.6667
$ws.Range("K86").Value = 2280.4546
$ws.Range("L86").Value = 2915.6667
$ws.Range("M86").Value = -1157.4546
$ws.Range("N86").Value = -5161.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2378.1794
$ws.Range("I89").Value = 2280.4546
$ws.Range("J89").Value = 2915.6667
$ws.Range("K89").Value = 11402.273
$ws.Range("L89").Value = 14578.3335
$ws.Range("M89").Value = -5786.273000000001
$ws.Range("N89").Value = -25810.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2967.75
$ws.Range("I94").Value = 1520.6875
$ws.Range("J94").Value = 8756
$ws.Range("K94").Value = 1520.6875
$ws.Range("L94").Value = 8756
$ws.Range("M94").Value = -1069.6875
$ws.Range("N94").Value = -9658

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3565.7273
$ws.Range("I99").Value = 3407.875
$ws.Range("J99").Value = 3986.6667
$ws.Range("K99").Value = 3407.875
$ws.Range("L99").Value = 3986.6667
$ws.Range("M99").Value = -1909.875
$ws.Range("N99").Value = -6982.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4526.9287
$ws.Range("I107").Value = 4332.5
$ws.Range("J107").Value = 4915.7856
$ws.Range("K107").Value = 4332.5
$ws.Range("L107").Value = 4915.7856
$ws.Range("M107").Value = -2412.5
$ws.Range("N107").Value = -8755.785599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2541.9092
$ws.Range("I134").Value = 2556.1
$ws.Range("K134").Value = 7668.299999999999
$ws.Range("M134").Value = -5133.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 5204
$ws.Range("I17").Value = 408
$ws.Range("K17").Value = 408
$ws.Range("M17").Value = -234

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1506.375
$ws.Range("J31").Value = 1195.5
$ws.Range("L31").Value = 1195.5
$ws.Range("N31").Value = -1785.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1506.375
$ws.Range("J34").Value = 1195.5
$ws.Range("L34").Value = 1195.5
$ws.Range("N34").Value = -1599.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 41670316
$ws.Range("I86").Value = 62503504
$ws.Range("K86").Value = 62503504
$ws.Range("M86").Value = -62502381

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 41670316
$ws.Range("I89").Value = 62503504
$ws.Range("K89").Value = 312517520
$ws.Range("M89").Value = -312511904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1482.3
$ws.Range("I94").Value = 1454.909
$ws.Range("K94").Value = 1454.909
$ws.Range("M94").Value = -1003.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 8000858
$ws.Range("I107").Value = 36000400
$ws.Range("K107").Value = 36000400
$ws.Range("M107").Value = -35998480

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3271.75
$ws.Range("I132").Value = 3295.6667
$ws.Range("K132").Value = 9887.000100000001
$ws.Range("M132").Value = -7357.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 12932.667
$ws.Range("J106").Value = 13499
$ws.Range("L106").Value = 40497
$ws.Range("N106").Value = -42389

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1766.1875
$ws.Range("J131").Value = 1679.375
$ws.Range("L131").Value = 5038.125
$ws.Range("N131").Value = -15118.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 620
$ws.Range("I97").Value = 572.7692
$ws.Range("K97").Value = 572.7692
$ws.Range("M97").Value = -76.76919999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3500.4285
$ws.Range("I102").Value = 3126.25
$ws.Range("K102").Value = 3126.25
$ws.Range("M102").Value = -1504.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2208.1428
$ws.Range("I132").Value = 1215.2424
$ws.Range("J132").Value = 4256
$ws.Range("K132").Value = 3645.7272
$ws.Range("L132").Value = 12768
$ws.Range("M132").Value = -1115.7272
$ws.Range("N132").Value = -17828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9260262
$ws.Range("I46").Value = 11905631
$ws.Range("J46").Value = 1469.5
$ws.Range("K46").Value = 11905631
$ws.Range("L46").Value = 1469.5
$ws.Range("M46").Value = -11905443
$ws.Range("N46").Value = -1845.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 7364.2
$ws.Range("I93").Value = 6429
$ws.Range("J93").Value = 11105
$ws.Range("K93").Value = 6429
$ws.Range("L93").Value = 11105
$ws.Range("M93").Value = -5181
$ws.Range("N93").Value = -13601
